$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing rows 1142-1143 (O1142, R1142, R1143) ---
$ws.Cells.Item(1142, 15).Value = 1   # O1142: 0 -> 1
$ws.Cells.Item(1142, 18).Value = 0   # R1142: inlineStr(blank) -> 0
$ws.Cells.Item(1143, 18).Value = 0   # R1143: inlineStr(blank) -> 0

# --- Append new weekly rows 1144-1156 ---
# Row 1144
$ws.Cells.Item(1144, 1).Value = 45474
$ws.Cells.Item(1144, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1144, 2).Value = 1596.900024414062
$ws.Cells.Item(1144, 3).Value = 1608.949951171875
$ws.Cells.Item(1144, 4).Value = 1562.050048828125
$ws.Cells.Item(1144, 5).Value = 1592.75
$ws.Cells.Item(1144, 6).Value = 1586.875
$ws.Cells.Item(1144, 7).Value = 978670
$ws.Cells.Item(1144, 8).Value = 2024
$ws.Cells.Item(1144, 9).Value = 7
$ws.Cells.Item(1144, 10).Value = 1
$ws.Cells.Item(1144, 11).Value = 0
$ws.Cells.Item(1144, 12).Value = 0
$ws.Cells.Item(1144, 13).Value = 0
$ws.Cells.Item(1144, 14).Value = 27
$ws.Cells.Item(1144, 15).Value = 0
$ws.Cells.Item(1144, 16).Value = 0
$ws.Cells.Item(1144, 17).Value = 0

# Row 1145
$ws.Cells.Item(1145, 1).Value = 45481
$ws.Cells.Item(1145, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1145, 2).Value = 1599
$ws.Cells.Item(1145, 3).Value = 1638.800048828125
$ws.Cells.Item(1145, 4).Value = 1571.5
$ws.Cells.Item(1145, 5).Value = 1602
$ws.Cells.Item(1145, 6).Value = 1596.090942382812
$ws.Cells.Item(1145, 7).Value = 1815882
$ws.Cells.Item(1145, 8).Value = 2024
$ws.Cells.Item(1145, 9).Value = 7
$ws.Cells.Item(1145, 10).Value = 8
$ws.Cells.Item(1145, 11).Value = 0
$ws.Cells.Item(1145, 12).Value = 0
$ws.Cells.Item(1145, 13).Value = 0
$ws.Cells.Item(1145, 14).Value = 28
$ws.Cells.Item(1145, 15).Value = 0
$ws.Cells.Item(1145, 16).Value = 0
$ws.Cells.Item(1145, 17).Value = 0

# Row 1146
$ws.Cells.Item(1146, 1).Value = 45488
$ws.Cells.Item(1146, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1146, 2).Value = 1602.099975585938
$ws.Cells.Item(1146, 3).Value = 1622.650024414062
$ws.Cells.Item(1146, 4).Value = 1558.199951171875
$ws.Cells.Item(1146, 5).Value = 1595.400024414062
$ws.Cells.Item(1146, 6).Value = 1589.515258789062
$ws.Cells.Item(1146, 7).Value = 1559587
$ws.Cells.Item(1146, 8).Value = 2024
$ws.Cells.Item(1146, 9).Value = 7
$ws.Cells.Item(1146, 10).Value = 15
$ws.Cells.Item(1146, 11).Value = 0
$ws.Cells.Item(1146, 12).Value = 0
$ws.Cells.Item(1146, 13).Value = 0
$ws.Cells.Item(1146, 14).Value = 29
$ws.Cells.Item(1146, 15).Value = 0
$ws.Cells.Item(1146, 16).Value = 0
$ws.Cells.Item(1146, 17).Value = 0

# Row 1147
$ws.Cells.Item(1147, 1).Value = 45495
$ws.Cells.Item(1147, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1147, 2).Value = 1595.400024414062
$ws.Cells.Item(1147, 3).Value = 1646.400024414062
$ws.Cells.Item(1147, 4).Value = 1499.699951171875
$ws.Cells.Item(1147, 5).Value = 1624.699951171875
$ws.Cells.Item(1147, 6).Value = 1618.707153320312
$ws.Cells.Item(1147, 7).Value = 3668893
$ws.Cells.Item(1147, 8).Value = 2024
$ws.Cells.Item(1147, 9).Value = 7
$ws.Cells.Item(1147, 10).Value = 22
$ws.Cells.Item(1147, 11).Value = 0
$ws.Cells.Item(1147, 12).Value = 0
$ws.Cells.Item(1147, 13).Value = 0
$ws.Cells.Item(1147, 14).Value = 30
$ws.Cells.Item(1147, 15).Value = 0
$ws.Cells.Item(1147, 16).Value = 0
$ws.Cells.Item(1147, 17).Value = 0

# Row 1148
$ws.Cells.Item(1148, 1).Value = 45502
$ws.Cells.Item(1148, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1148, 2).Value = 1632.199951171875
$ws.Cells.Item(1148, 3).Value = 1694.400024414062
$ws.Cells.Item(1148, 4).Value = 1605
$ws.Cells.Item(1148, 5).Value = 1638.949951171875
$ws.Cells.Item(1148, 6).Value = 1638.949951171875
$ws.Cells.Item(1148, 7).Value = 1817796
$ws.Cells.Item(1148, 8).Value = 2024
$ws.Cells.Item(1148, 9).Value = 7
$ws.Cells.Item(1148, 10).Value = 29
$ws.Cells.Item(1148, 11).Value = 0
$ws.Cells.Item(1148, 12).Value = 0
$ws.Cells.Item(1148, 13).Value = 0
$ws.Cells.Item(1148, 14).Value = 31
$ws.Cells.Item(1148, 15).Value = 0
$ws.Cells.Item(1148, 16).Value = 0
$ws.Cells.Item(1148, 17).Value = 2

# Row 1149
$ws.Cells.Item(1149, 1).Value = 45509
$ws.Cells.Item(1149, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1149, 2).Value = 1610
$ws.Cells.Item(1149, 3).Value = 1704
$ws.Cells.Item(1149, 4).Value = 1585.550048828125
$ws.Cells.Item(1149, 5).Value = 1693.699951171875
$ws.Cells.Item(1149, 6).Value = 1693.699951171875
$ws.Cells.Item(1149, 7).Value = 4731253
$ws.Cells.Item(1149, 8).Value = 2024
$ws.Cells.Item(1149, 9).Value = 8
$ws.Cells.Item(1149, 10).Value = 5
$ws.Cells.Item(1149, 11).Value = 0
$ws.Cells.Item(1149, 12).Value = 0
$ws.Cells.Item(1149, 13).Value = 0
$ws.Cells.Item(1149, 14).Value = 32
$ws.Cells.Item(1149, 15).Value = 0
$ws.Cells.Item(1149, 16).Value = 0
$ws.Cells.Item(1149, 17).Value = 0

# Row 1150
$ws.Cells.Item(1150, 1).Value = 45516
$ws.Cells.Item(1150, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1150, 2).Value = 1709
$ws.Cells.Item(1150, 3).Value = 1780
$ws.Cells.Item(1150, 4).Value = 1673.550048828125
$ws.Cells.Item(1150, 5).Value = 1759.75
$ws.Cells.Item(1150, 6).Value = 1759.75
$ws.Cells.Item(1150, 7).Value = 2484808
$ws.Cells.Item(1150, 8).Value = 2024
$ws.Cells.Item(1150, 9).Value = 8
$ws.Cells.Item(1150, 10).Value = 12
$ws.Cells.Item(1150, 11).Value = 0
$ws.Cells.Item(1150, 12).Value = 0
$ws.Cells.Item(1150, 13).Value = 0
$ws.Cells.Item(1150, 14).Value = 33
$ws.Cells.Item(1150, 15).Value = 0
$ws.Cells.Item(1150, 16).Value = 0
$ws.Cells.Item(1150, 17).Value = 0

# Row 1151
$ws.Cells.Item(1151, 1).Value = 45523
$ws.Cells.Item(1151, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1151, 2).Value = 1773.949951171875
$ws.Cells.Item(1151, 3).Value = 1789
$ws.Cells.Item(1151, 4).Value = 1713.349975585938
$ws.Cells.Item(1151, 5).Value = 1748.5
$ws.Cells.Item(1151, 6).Value = 1748.5
$ws.Cells.Item(1151, 7).Value = 2019495
$ws.Cells.Item(1151, 8).Value = 2024
$ws.Cells.Item(1151, 9).Value = 8
$ws.Cells.Item(1151, 10).Value = 19
$ws.Cells.Item(1151, 11).Value = 0
$ws.Cells.Item(1151, 12).Value = 0
$ws.Cells.Item(1151, 13).Value = 0
$ws.Cells.Item(1151, 14).Value = 34
$ws.Cells.Item(1151, 15).Value = 1
$ws.Cells.Item(1151, 16).Value = 0
$ws.Cells.Item(1151, 17).Value = 0

# Row 1152
$ws.Cells.Item(1152, 1).Value = 45530
$ws.Cells.Item(1152, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1152, 2).Value = 1748.5
$ws.Cells.Item(1152, 3).Value = 1773.75
$ws.Cells.Item(1152, 4).Value = 1706.050048828125
$ws.Cells.Item(1152, 5).Value = 1758.25
$ws.Cells.Item(1152, 6).Value = 1758.25
$ws.Cells.Item(1152, 7).Value = 1593884
$ws.Cells.Item(1152, 8).Value = 2024
$ws.Cells.Item(1152, 9).Value = 8
$ws.Cells.Item(1152, 10).Value = 26
$ws.Cells.Item(1152, 11).Value = 0
$ws.Cells.Item(1152, 12).Value = 0
$ws.Cells.Item(1152, 13).Value = 0
$ws.Cells.Item(1152, 14).Value = 35
$ws.Cells.Item(1152, 15).Value = 0
$ws.Cells.Item(1152, 16).Value = 0
$ws.Cells.Item(1152, 17).Value = 0

# Row 1153
$ws.Cells.Item(1153, 1).Value = 45537
$ws.Cells.Item(1153, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1153, 2).Value = 1758.25
$ws.Cells.Item(1153, 3).Value = 1773.550048828125
$ws.Cells.Item(1153, 4).Value = 1681
$ws.Cells.Item(1153, 5).Value = 1691.150024414062
$ws.Cells.Item(1153, 6).Value = 1691.150024414062
$ws.Cells.Item(1153, 7).Value = 1452891
$ws.Cells.Item(1153, 8).Value = 2024
$ws.Cells.Item(1153, 9).Value = 9
$ws.Cells.Item(1153, 10).Value = 2
$ws.Cells.Item(1153, 11).Value = 0
$ws.Cells.Item(1153, 12).Value = 0
$ws.Cells.Item(1153, 13).Value = 0
$ws.Cells.Item(1153, 14).Value = 36
$ws.Cells.Item(1153, 15).Value = 0
$ws.Cells.Item(1153, 16).Value = 0
$ws.Cells.Item(1153, 17).Value = 0

# Row 1154
$ws.Cells.Item(1154, 1).Value = 45544
$ws.Cells.Item(1154, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1154, 2).Value = 1688
$ws.Cells.Item(1154, 3).Value = 1728.25
$ws.Cells.Item(1154, 4).Value = 1667.050048828125
$ws.Cells.Item(1154, 5).Value = 1694.400024414062
$ws.Cells.Item(1154, 6).Value = 1694.400024414062
$ws.Cells.Item(1154, 7).Value = 1318926
$ws.Cells.Item(1154, 8).Value = 2024
$ws.Cells.Item(1154, 9).Value = 9
$ws.Cells.Item(1154, 10).Value = 9
$ws.Cells.Item(1154, 11).Value = 0
$ws.Cells.Item(1154, 12).Value = 0
$ws.Cells.Item(1154, 13).Value = 0
$ws.Cells.Item(1154, 14).Value = 37
$ws.Cells.Item(1154, 15).Value = 0
$ws.Cells.Item(1154, 16).Value = 0
$ws.Cells.Item(1154, 17).Value = 0

# Row 1155
$ws.Cells.Item(1155, 1).Value = 45551
$ws.Cells.Item(1155, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1155, 2).Value = 1694.400024414062
$ws.Cells.Item(1155, 3).Value = 1743.25
$ws.Cells.Item(1155, 4).Value = 1606.900024414062
$ws.Cells.Item(1155, 5).Value = 1659.150024414062
$ws.Cells.Item(1155, 6).Value = 1659.150024414062
$ws.Cells.Item(1155, 7).Value = 2535435
$ws.Cells.Item(1155, 8).Value = 2024
$ws.Cells.Item(1155, 9).Value = 9
$ws.Cells.Item(1155, 10).Value = 16
$ws.Cells.Item(1155, 11).Value = 0
$ws.Cells.Item(1155, 12).Value = 0
$ws.Cells.Item(1155, 13).Value = 0
$ws.Cells.Item(1155, 14).Value = 38
$ws.Cells.Item(1155, 15).Value = 0
$ws.Cells.Item(1155, 16).Value = 0
$ws.Cells.Item(1155, 17).Value = 0

# Row 1156
$ws.Cells.Item(1156, 1).Value = 45558
$ws.Cells.Item(1156, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1156, 2).Value = 1667.900024414062
$ws.Cells.Item(1156, 3).Value = 1686
$ws.Cells.Item(1156, 4).Value = 1612
$ws.Cells.Item(1156, 5).Value = 1657
$ws.Cells.Item(1156, 6).Value = 1657
$ws.Cells.Item(1156, 7).Value = 2249211
$ws.Cells.Item(1156, 8).Value = 2024
$ws.Cells.Item(1156, 9).Value = 9
$ws.Cells.Item(1156, 10).Value = 23
$ws.Cells.Item(1156, 11).Value = 0
$ws.Cells.Item(1156, 12).Value = 0
$ws.Cells.Item(1156, 13).Value = 0
$ws.Cells.Item(1156, 14).Value = 39
$ws.Cells.Item(1156, 15).Value = 0
$ws.Cells.Item(1156, 16).Value = 0
$ws.Cells.Item(1156, 17).Value = 0
